$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "_spritePath" (column H) values for each item row.
# Cells are written in this specific order so that new shared-string
# entries are appended to xl/sharedStrings.xml in the same sequence
# as in the target workbook.
$ws.Range("H3").Value  = "Items/CandyTree_Item"
$ws.Range("H4").Value  = "Items/Man-EatingFlower_Item"
$ws.Range("H6").Value  = "Items/Grape_Item"
$ws.Range("H5").Value  = "Items/Herb_Item"
$ws.Range("H7").Value  = "Items/Mushroom_Item"
$ws.Range("H9").Value  = "Items/Turtle_Item"
$ws.Range("H10").Value = "Items/Grasshopper_Item"
$ws.Range("H11").Value = "Items/Firelizard_Item"
$ws.Range("H12").Value = "Items/Corgi_Item"
$ws.Range("H13").Value = "Items/Kirby_Item"
$ws.Range("H18").Value = "CombineItems/Emerald_Potion"
$ws.Range("H19").Value = "CombineItems/Ruby_Potion"
$ws.Range("H20").Value = "CombineItems/Star_Candy"
$ws.Range("H21").Value = "CombineItems/Pink_Protain_Shake"
$ws.Range("H22").Value = "CombineItems/Ostin-Corgi_Jam"
$ws.Range("H23").Value = "CombineItems/Mint_Stock"
$ws.Range("H8").Value  = "Items/FlySalamander_Item"
$ws.Range("H14").Value = "Items/LittleGlass_Item"
$ws.Range("H15").Value = "Items/BigGlass_Item"
$ws.Range("H16").Value = "Items/JamGlass_Item"
$ws.Range("H17").Value = "Items/CandyWrap_Item"

# Resize column H to fit the new (longer) sprite-path strings.
$ws.Columns.Item(8).EntireColumn.AutoFit()

# Match the saved selection state from the authored workbook.
$null = $ws.Range("H18").Select()
